# Updates the cryptos worksheet with the latest scraped price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.318.43'
$ws.Range('E2').Value = '  -3.00%  '
$ws.Range('D3').Value = '1.975.50'
$ws.Range('E3').Value = '  -3.77%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = "'245.57"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.62%  '
$ws.Range('D6').Value = "'0.624"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.77%  '
$ws.Range('D7').Value = "'58.73"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -11.16%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = "'0.372"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.88%  '
$ws.Range('D10').Value = "'56.52"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.04%  '
$ws.Range('D11').Value = "'0.0838"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.44%  '
$ws.Range('E12').Value = '  -0.68%  '
$ws.Range('D13').Value = "'22.96"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.42%  '
$ws.Range('D14').Value = "'0.858"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.99%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = "'13.90"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.34%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '2.259.90'
$ws.Range('E16').Value = '  -4.03%  '
$ws.Range('D17').Value = "'5.43"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.00%  '
$ws.Range('D18').Value = '1.973.44'
$ws.Range('E18').Value = '  -4.42%  '
$ws.Range('D19').Value = '36.207.47'
$ws.Range('E19').Value = '  -3.17%  '
$ws.Range('D20').Value = "'70.45"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.67%  '
$ws.Range('D21').Value = '0.0₃0878'
$ws.Range('E21').Value = '  -1.53%  '
$ws.Range('D22').Value = "'5.26"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.77%  '
$ws.Range('D23').Value = "'233.83"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.22%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').Value = "'2.51"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.89%  '
$ws.Range('D26').Value = "'2.31"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.29%  '
$ws.Range('D27').Value = "'9.81"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.70%  '
$ws.Range('D28').Value = "'163.61"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('D29').Value = "'0.134"
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Value = "'19.73"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.28%  '
$ws.Range('E31').Value = '  -3.13%  '
$ws.Range('D32').Value = "'1.17"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.52%  '
$ws.Range('D33').Value = "'4.85"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.57%  '
$ws.Range('D34').Value = "'0.0654"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.85%  '
$ws.Range('D35').Value = "'4.43"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.24%  '
$ws.Range('B36').Value = 'THORChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D36').Value = "'6.13"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.65%  '
$ws.Range('B37').Value = 'BinanceUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D37').Value = "'1.00"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('E38').Value = '  -1.31%  '
$ws.Range('D39').Value = "'2.22"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -9.16%  '
$ws.Range('D40').Value = "'2.92"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.01%  '
$ws.Range('D41').Value = "'1.22"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.08%  '
$ws.Range('D42').Value = "'0.0959"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.41%  '
$ws.Range('D43').Value = "'2.88"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.30%  '
$ws.Range('D44').Value = "'0.0213"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.64%  '
$ws.Range('E45').Value = '  -7.49%  '
$ws.Range('D46').Value = "'16.04"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -10.13%  '
$ws.Range('D47').Value = "'91.30"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.20%  '
$ws.Range('D48').Value = '1.363.51'
$ws.Range('E48').Value = '  -3.72%  '
$ws.Range('D49').Value = "'7.40"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.59%  '
$ws.Range('E50').Value = '  -4.33%  '
$ws.Range('D51').Value = "'44.96"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.23%  '
